$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bulk-update column C (rows 2-322): "Förändrad" date 45184 -> 45186.
$ws.Range("C2:C322").Value2 = 45186

# 2. Add a friendly display-name second argument to every HYPERLINK() formula
#    in columns S,T,U,V,W,X,Y (rows 2-17 are the only ones that currently
#    carry these formulas). The display name is the row's "Beteckning"
#    (column A) value, which is also the filename stem used in each URL.
$linkCols = @("S","T","U","V","W","X","Y")
for ($r = 2; $r -le 17; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $r)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -match '^(=HYPERLINK\("[^"]+")\)$') {
                $cell.Formula = $matches[1] + ', "' + $beteckning + '")'
            }
        }
    }
}

# 3. Row 322 now keeps the sheet's standard (default) row height explicitly.
$ws.Rows(322).RowHeight = 15

# 4. Append the new row (323) with the latest case.
$newRow = 323
$ws.Cells.Item($newRow, 1).Value2 = "A 43485-2023"
$ws.Cells.Item($newRow, 2).Value2 = 45184
$ws.Cells.Item($newRow, 3).Value2 = 45186
$ws.Cells.Item($newRow, 4).Value2 = "JÖNKÖPINGS LÄN"
$ws.Cells.Item($newRow, 5).Value2 = "ANEBY"
$ws.Cells.Item($newRow, 7).Value2 = 3.4
for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item($newRow, $c).Value2 = 0
}
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"
$ws.Range("R" + $newRow).WrapText = $true
